$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) — rename existing handback file and add a new one
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Duplicate row 2 down into row 3 (keeps cell styles / blank-column layout)
$wsOverview.Rows.Item(2).Copy()
$wsOverview.Rows.Item(3).Insert()

# Update row 2 (36db075a... -> 2fba2052...)
$wsOverview.Range("A2").Value = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.md"
$wsOverview.Range("B2").Value = "e2e\2fba2052-50bf-4790-8c8c-2e07826b5dc7.md"
$wsOverview.Range("G2").Value = "2016-08-24 03:00:52"

# Replace the hyperlink on B2 so it points at the renamed file
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4ad9baa6fc4ad34c030d70560d49ba1811b1e65/e2e/2fba2052-50bf-4790-8c8c-2e07826b5dc7.md", "", "", "e2e\2fba2052-50bf-4790-8c8c-2e07826b5dc7.md") | Out-Null

# Fill in row 3 (new 68e22ea2... handback entry)
$wsOverview.Range("A3").Value = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md"
$wsOverview.Range("B3").Value = "e2e\68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-24 03:00:52"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4ad9baa6fc4ad34c030d70560d49ba1811b1e65/e2e/68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md", "", "", "e2e\68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md") | Out-Null

# Grow the "Overview" table to include the new row
$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(2).Copy()
$wsZh.Rows.Item(3).Insert()

# Update row 2 (36db075a... -> 2fba2052..., new xliff hash + dates)
$wsZh.Range("A2").Value = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.md"
$wsZh.Range("G2").Value = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.ef624036b7a2d9357573e67962c90ab2b542cc3e.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-24 03:00:47"
$wsZh.Range("I2").Value = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.md"
$wsZh.Range("J2").Value = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.ef624036b7a2d9357573e67962c90ab2b542cc3e.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-24 03:01:17"

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4ad9baa6fc4ad34c030d70560d49ba1811b1e65/e2e/2fba2052-50bf-4790-8c8c-2e07826b5dc7.md", "", "", "2fba2052-50bf-4790-8c8c-2e07826b5dc7.md") | Out-Null
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/939bfe627b4798f4f78bf65e3b2bbbce9de57613/e2e/2fba2052-50bf-4790-8c8c-2e07826b5dc7.md", "", "", "2fba2052-50bf-4790-8c8c-2e07826b5dc7.md") | Out-Null

# Fill in row 3 (new 68e22ea2... handback entry)
$wsZh.Range("A3").Value = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.1e28be1cea68a5b58d1cc99c5bc006ea64757a9e.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-24 03:00:47"
$wsZh.Range("I3").Value = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md"
$wsZh.Range("J3").Value = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.1e28be1cea68a5b58d1cc99c5bc006ea64757a9e.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-24 03:01:17"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4ad9baa6fc4ad34c030d70560d49ba1811b1e65/e2e/68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md", "", "", "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/939bfe627b4798f4f78bf65e3b2bbbce9de57613/e2e/68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md", "", "", "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md") | Out-Null

$wsZh.ListObjects.Item(1).Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(2).Copy()
$wsDe.Rows.Item(3).Insert()

# Update row 2 (36db075a... -> 2fba2052..., new xliff hash + dates)
$wsDe.Range("A2").Value = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.md"
$wsDe.Range("G2").Value = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.ef624036b7a2d9357573e67962c90ab2b542cc3e.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-24 03:00:52"
$wsDe.Range("I2").Value = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.md"
$wsDe.Range("J2").Value = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.ef624036b7a2d9357573e67962c90ab2b542cc3e.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-24 03:01:24"

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4ad9baa6fc4ad34c030d70560d49ba1811b1e65/e2e/2fba2052-50bf-4790-8c8c-2e07826b5dc7.md", "", "", "2fba2052-50bf-4790-8c8c-2e07826b5dc7.md") | Out-Null
$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f050b4d792a3a35ec8d90d04946c19cde9c5a49b/e2e/2fba2052-50bf-4790-8c8c-2e07826b5dc7.md", "", "", "2fba2052-50bf-4790-8c8c-2e07826b5dc7.md") | Out-Null

# Fill in row 3 (new 68e22ea2... handback entry)
$wsDe.Range("A3").Value = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.1e28be1cea68a5b58d1cc99c5bc006ea64757a9e.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-24 03:00:52"
$wsDe.Range("I3").Value = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md"
$wsDe.Range("J3").Value = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.1e28be1cea68a5b58d1cc99c5bc006ea64757a9e.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-24 03:01:24"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4ad9baa6fc4ad34c030d70560d49ba1811b1e65/e2e/68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md", "", "", "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f050b4d792a3a35ec8d90d04946c19cde9c5a49b/e2e/68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md", "", "", "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md") | Out-Null

$wsDe.ListObjects.Item(1).Resize($wsDe.Range("A1:P3"))
